$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet -------------------------------------------------
$ws.Name = "Login Page"

# --- Update / add cell values -----------------------------------------
# Row 2: describe the Login Page functionality
$ws.Range("C2").Value = "User should be able to login or create user"

# Row 3: Signup -> Signup Page, plus new test data/status/tester info
$ws.Range("C3").Value = "Signup Page"
$ws.Range("D3").Value = "User correct data"
$ws.Range("E3").Value = "Failed"
$ws.Range("F3").Value = "Pavan"

# Row 1 header: add a new "Current Status" column in H
$ws.Range("H1").Value = "Current Status"

# Rows 7 & 10: mark the fixed issues in the new "Current Status" column
$ws.Range("H7").Value = "Fixed"
$ws.Range("H10").Value = "Fixed"

# --- Match header formatting (F1:G1 become bold header style like A1:E1)
$ws.Range("A1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)  # xlPasteFormats

# --- Selection ----------------------------------------------------------
$ws.Range("H10").Select() | Out-Null

# --- AutoFilter -----------------------------------------------------------
[void]$ws.Range("A1:G10").AutoFilter()

# --- Hidden _FilterDatabase defined name (what AutoFilter implies in Excel)
$n = $ws.Names.Add("_xlnm._FilterDatabase", "='Login Page'!`$A`$1:`$G`$10")
$n.Visible = $false
